$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-33: columns B (Mid X), C (Mid Y), E (Rot) ---
$updates = @(
    @{ Row = 2;  B = 15;   C = 87;   E = 270 },
    @{ Row = 3;  B = 15;   C = 77;   E = 270 },
    @{ Row = 4;  B = 15;   C = 67;   E = 270 },
    @{ Row = 5;  B = 15;   C = 57;   E = 270 },
    @{ Row = 6;  B = 15;   C = 47;   E = 270 },
    @{ Row = 7;  B = 15;   C = 37;   E = 270 },
    @{ Row = 8;  B = 15;   C = 27;   E = 270 },
    @{ Row = 9;  B = 15;   C = 17;   E = 270 },
    @{ Row = 10; B = 22.5; C = 87;   E = 270 },
    @{ Row = 11; B = 22.5; C = 77;   E = 270 },
    @{ Row = 12; B = 19.5; C = 67;   E = 0 },
    @{ Row = 13; B = 19.5; C = 57;   E = 0 },
    @{ Row = 14; B = 19.5; C = 47;   E = 0 },
    @{ Row = 15; B = 19.5; C = 37;   E = 0 },
    @{ Row = 16; B = 19.5; C = 27;   E = 0 },
    @{ Row = 17; B = 19.5; C = 17;   E = 0 },
    @{ Row = 18; B = 31.5; C = 84.5; E = 0 },
    @{ Row = 19; B = 31.5; C = 74.5; E = 0 },
    @{ Row = 20; B = 31.5; C = 89.5; E = 0 },
    @{ Row = 21; B = 31.5; C = 79.5; E = 0 },
    @{ Row = 22; B = 23.5; C = 66;   E = 0 },
    @{ Row = 23; B = 23.5; C = 56;   E = 0 },
    @{ Row = 24; B = 23.5; C = 68;   E = 0 },
    @{ Row = 25; B = 23.5; C = 58;   E = 0 },
    @{ Row = 26; B = 23.5; C = 46;   E = 0 },
    @{ Row = 27; B = 23.5; C = 36;   E = 0 },
    @{ Row = 28; B = 23.5; C = 48;   E = 0 },
    @{ Row = 29; B = 23.5; C = 38;   E = 0 },
    @{ Row = 30; B = 23.5; C = 26;   E = 0 },
    @{ Row = 31; B = 23.5; C = 16;   E = 0 },
    @{ Row = 32; B = 23.5; C = 28;   E = 0 },
    @{ Row = 33; B = 23.5; C = 18;   E = 0 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 5).Value = $u.E
}

# --- Append new designator rows 34-46 (R17..R29) ---
$newRows = @(
    @{ Row = 34; Designator = "R17"; B = 59.29; C = 87;   D = "top"; E = 0 },
    @{ Row = 35; Designator = "R18"; B = 59.29; C = 84;   D = "top"; E = 0 },
    @{ Row = 36; Designator = "R19"; B = 70.29; C = 87;   D = "top"; E = 0 },
    @{ Row = 37; Designator = "R20"; B = 70.29; C = 84;   D = "top"; E = 0 },
    @{ Row = 38; Designator = "R21"; B = 83.81; C = 14.5; D = "top"; E = 0 },
    @{ Row = 39; Designator = "R22"; B = 83.81; C = 89.5; D = "top"; E = 0 },
    @{ Row = 40; Designator = "R23"; B = 83.81; C = 79.5; D = "top"; E = 0 },
    @{ Row = 41; Designator = "R24"; B = 83.81; C = 69.5; D = "top"; E = 0 },
    @{ Row = 42; Designator = "R25"; B = 83.81; C = 59.5; D = "top"; E = 0 },
    @{ Row = 43; Designator = "R26"; B = 83.81; C = 49.5; D = "top"; E = 0 },
    @{ Row = 44; Designator = "R27"; B = 83.81; C = 39.5; D = "top"; E = 0 },
    @{ Row = 45; Designator = "R28"; B = 83.81; C = 29.5; D = "top"; E = 0 },
    @{ Row = 46; Designator = "R29"; B = 83.81; C = 19.5; D = "top"; E = 0 }
)

foreach ($n in $newRows) {
    $r = $n.Row
    $ws.Cells.Item($r, 1).Value = $n.Designator
    $ws.Cells.Item($r, 2).Value = $n.B
    $ws.Cells.Item($r, 3).Value = $n.C
    $ws.Cells.Item($r, 4).Value = $n.D
    $ws.Cells.Item($r, 5).Value = $n.E
}

# --- Update the selection to mirror the authored state ---
$ws.Range("A40").Select()
